# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gets a new blank column inserted before
# column N (shifting the old "Late" / "heading" / "Outstanding" columns
# one slot to the right: N->O, O->P, P->Q), and the active selection
# moves to R7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Capture column M's width so the freshly-inserted column N inherits the
# same explicit width (matching the neighbouring "In Advance" column)
# instead of picking up the width that used to belong to old column N.
$mWidth = $ws.Columns("M").ColumnWidth

# Insert a new blank column before N; everything at/after N shifts right.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Restore the active cell selection to where it ended up after the edit.
$ws.Range("R7").Select() | Out-Null
